# "collection and quality reports v1"
#
# Row 2 on Sheet1 is repurposed from a stray test login (shritej.m / 123 /
# the Panvel offline-payment URL) to the real Baramati MC login row, and the
# Url cell (C2) becomes a live hyperlink to that page. Finally the sheet
# selection is left on A2:B2 (the cells that were just edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the hyperlink before touching the cell's value/content. Excel
# stamps the target cell with its "Hyperlink" cell style as a side effect of
# Hyperlinks.Add; nudging a harmless formatting property afterwards (it's
# already centered, same as the neighbouring C3 url cell) makes Excel
# re-resolve that style back onto the existing shared "Hyperlink" style (s=7)
# instead of leaving a redundant duplicate applied to the cell.
$ws.Hyperlinks.Add($ws.Range("C2"), "http://testbaramatimc.ptaxcollection.com:8080/Pages/OfflinePayment.aspx") | Out-Null

$ws.Range("A2").Value = "jagdish.d"
$ws.Range("B2").Value = "User@12345"
$ws.Range("C2").Value = "http://testbaramatimc.ptaxcollection.com:8080/Pages/OfflinePayment.aspx"

$ws.Range("C2").HorizontalAlignment = $ws.Range("C3").HorizontalAlignment

# Leave the selection on the cells that were just edited.
$ws.Range("A2:B2").Select() | Out-Null
